$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.41485066666667
$ws.Range("H2").Value = 52.24455200000001
$ws.Range("I2").Value = 0.1047285618770465
$ws.Range("J2").Value = 0.1047285618770465
$ws.Range("M2").Value = 4.717738333333333
$ws.Range("N2").Value = 14.153215
$ws.Range("O2").Value = 0.2002263444295212
$ws.Range("P2").Value = 0.2002263444295212
$ws.Range("Q2").Value = 82.15870855940891
$ws.Range("R2").Value = 739.42837703468
$ws.Range("S2").Value = 0.02096941710200194
$ws.Range("T2").Value = 0.02096941710200194
$ws.Range("G3").Value = 17.41485066666667
$ws.Range("H3").Value = 52.24455200000001
$ws.Range("I3").Value = 0.1047285618770465
$ws.Range("J3").Value = 0.1047285618770465
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.1305610278731266
$ws.Range("P3").Value = 0.1305610278731266
$ws.Range("Q3").Value = 53.5729974435049
$ws.Range("R3").Value = 482.1569769915441
$ws.Range("S3").Value = 0.01367346868634154
$ws.Range("T3").Value = 0.01367346868634154
$ws.Range("G4").Value = 17.41485066666667
$ws.Range("H4").Value = 52.24455200000001
$ws.Range("I4").Value = 0.1047285618770465
$ws.Range("J4").Value = 0.1047285618770465
$ws.Range("M4").Value = 0.6908423333333333
$ws.Range("N4").Value = 2.072527
$ws.Range("O4").Value = 0.02932015834857891
$ws.Range("P4").Value = 0.02932015834857891
$ws.Range("Q4").Value = 12.03091606921156
$ws.Range("R4").Value = 108.278244622904
$ws.Range("S4").Value = 0.003070658017853949
$ws.Range("T4").Value = 0.003070658017853949
$ws.Range("G5").Value = 17.41485066666667
$ws.Range("H5").Value = 52.24455200000001
$ws.Range("I5").Value = 0.1047285618770465
$ws.Range("J5").Value = 0.1047285618770465
$ws.Range("M5").Value = 15.077163
$ws.Range("N5").Value = 45.231489
$ws.Range("O5").Value = 0.6398924693487733
$ws.Range("P5").Value = 0.6398924693487733
$ws.Range("Q5").Value = 262.566542121992
$ws.Range("R5").Value = 2363.098879097928
$ws.Range("S5").Value = 0.0670150180708491
$ws.Range("T5").Value = 0.0670150180708491
$ws.Range("I6").Value = 0.1785014126970782
$ws.Range("J6").Value = 0.1785014126970782
$ws.Range("M6").Value = 4.717738333333333
$ws.Range("N6").Value = 14.153215
$ws.Range("O6").Value = 0.2002263444295212
$ws.Range("P6").Value = 0.2002263444295212
$ws.Range("Q6").Value = 140.0329125156856
$ws.Range("R6").Value = 1260.29621264117
$ws.Range("S6").Value = 0.03574068533984127
$ws.Range("T6").Value = 0.03574068533984129
$ws.Range("I7").Value = 0.1785014126970782
$ws.Range("J7").Value = 0.1785014126970782
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.1305610278731266
$ws.Range("P7").Value = 0.1305610278731266
$ws.Range("Q7").Value = 91.31086644070956
$ws.Range("R7").Value = 821.797797966386
$ws.Range("S7").Value = 0.0233053279185357
$ws.Range("T7").Value = 0.0233053279185357
$ws.Range("I8").Value = 0.1785014126970782
$ws.Range("J8").Value = 0.1785014126970782
$ws.Range("M8").Value = 0.6908423333333333
$ws.Range("N8").Value = 2.072527
$ws.Range("O8").Value = 0.02932015834857891
$ws.Range("P8").Value = 0.02932015834857891
$ws.Range("Q8").Value = 20.50572905713623
$ws.Range("R8").Value = 184.551561514226
$ws.Range("S8").Value = 0.005233689685723366
$ws.Range("T8").Value = 0.005233689685723367
$ws.Range("I9").Value = 0.1785014126970782
$ws.Range("J9").Value = 0.1785014126970782
$ws.Range("M9").Value = 15.077163
$ws.Range("N9").Value = 45.231489
$ws.Range("O9").Value = 0.6398924693487733
$ws.Range("P9").Value = 0.6398924693487733
$ws.Range("Q9").Value = 447.523558575998
$ws.Range("R9").Value = 4027.712027183982
$ws.Range("S9").Value = 0.1142217097529778
$ws.Range("T9").Value = 0.1142217097529779
$ws.Range("G10").Value = 84.03051233333333
$ws.Range("H10").Value = 252.091537
$ws.Range("I10").Value = 0.5053385113032314
$ws.Range("J10").Value = 0.5053385113032314
$ws.Range("M10").Value = 4.717738333333333
$ws.Range("N10").Value = 14.153215
$ws.Range("O10").Value = 0.2002263444295212
$ws.Range("P10").Value = 0.2002263444295212
$ws.Range("Q10").Value = 396.4339692046061
$ws.Range("R10").Value = 3567.905722841455
$ws.Range("S10").Value = 0.1011820828177023
$ws.Range("T10").Value = 0.1011820828177023
$ws.Range("G11").Value = 84.03051233333333
$ws.Range("H11").Value = 252.091537
$ws.Range("I11").Value = 0.5053385113032314
$ws.Range("J11").Value = 0.5053385113032314
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.1305610278731266
$ws.Range("P11").Value = 0.1305610278731266
$ws.Range("Q11").Value = 258.5015805519821
$ws.Range("R11").Value = 2326.514224967839
$ws.Range("S11").Value = 0.06597751545962549
$ws.Range("T11").Value = 0.06597751545962549
$ws.Range("G12").Value = 84.03051233333333
$ws.Range("H12").Value = 252.091537
$ws.Range("I12").Value = 0.5053385113032314
$ws.Range("J12").Value = 0.5053385113032314
$ws.Range("M12").Value = 0.6908423333333333
$ws.Range("N12").Value = 2.072527
$ws.Range("O12").Value = 0.02932015834857891
$ws.Range("P12").Value = 0.02932015834857891
$ws.Range("Q12").Value = 58.05183521155544
$ws.Range("R12").Value = 522.4665169039989
$ws.Range("S12").Value = 0.01481660517104588
$ws.Range("T12").Value = 0.01481660517104588
$ws.Range("G13").Value = 84.03051233333333
$ws.Range("H13").Value = 252.091537
$ws.Range("I13").Value = 0.5053385113032314
$ws.Range("J13").Value = 0.5053385113032314
$ws.Range("M13").Value = 15.077163
$ws.Range("N13").Value = 45.231489
$ws.Range("O13").Value = 0.6398924693487733
$ws.Range("P13").Value = 0.6398924693487733
$ws.Range("Q13").Value = 1266.941731423177
$ws.Range("R13").Value = 11402.47558280859
$ws.Range("S13").Value = 0.3233623078548577
$ws.Range("T13").Value = 0.3233623078548577
$ws.Range("G14").Value = 35.158014
$ws.Range("H14").Value = 105.474042
$ws.Range("I14").Value = 0.2114315141226439
$ws.Range("J14").Value = 0.2114315141226439
$ws.Range("M14").Value = 4.717738333333333
$ws.Range("N14").Value = 14.153215
$ws.Range("O14").Value = 0.2002263444295212
$ws.Range("P14").Value = 0.2002263444295212
$ws.Range("Q14").Value = 165.86631037167
$ws.Range("R14").Value = 1492.79679334503
$ws.Range("S14").Value = 0.04233415916997566
$ws.Range("T14").Value = 0.04233415916997567
$ws.Range("G15").Value = 35.158014
$ws.Range("H15").Value = 105.474042
$ws.Range("I15").Value = 0.2114315141226439
$ws.Range("J15").Value = 0.2114315141226439
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.1305610278731266
$ws.Range("P15").Value = 0.1305610278731266
$ws.Range("Q15").Value = 108.155977343286
$ws.Range("R15").Value = 973.403796089574
$ws.Range("S15").Value = 0.02760471580862387
$ws.Range("T15").Value = 0.02760471580862387
$ws.Range("G16").Value = 35.158014
$ws.Range("H16").Value = 105.474042
$ws.Range("I16").Value = 0.2114315141226439
$ws.Range("J16").Value = 0.2114315141226439
$ws.Range("M16").Value = 0.6908423333333333
$ws.Range("N16").Value = 2.072527
$ws.Range("O16").Value = 0.02932015834857891
$ws.Range("P16").Value = 0.02932015834857891
$ws.Range("Q16").Value = 24.288644427126
$ws.Range("R16").Value = 218.597799844134
$ws.Range("S16").Value = 0.006199205473955717
$ws.Range("T16").Value = 0.006199205473955717
$ws.Range("G17").Value = 35.158014
$ws.Range("H17").Value = 105.474042
$ws.Range("I17").Value = 0.2114315141226439
$ws.Range("J17").Value = 0.2114315141226439
$ws.Range("M17").Value = 15.077163
$ws.Range("N17").Value = 45.231489
$ws.Range("O17").Value = 0.6398924693487733
$ws.Range("P17").Value = 0.6398924693487733
$ws.Range("Q17").Value = 530.083107834282
$ws.Range("R17").Value = 4770.747970508537
$ws.Range("S17").Value = 0.1352934336700886
$ws.Range("T17").Value = 0.1352934336700886
